# Timesheet sign-off: supervisor (Prakruti Sinha) fills in her name,
# initials and the date she signed off the timesheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# Supervisor Name field
$ws.Range("G6").Value = "Prakruti Sinha"

# Supervisor signature block: initials + sign-off date
$ws.Range("A27").Value = "P.S"
$ws.Range("D27").Value = 41682
$ws.Range("D25").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("H29").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
